# Scheduled-runner style market/profit data refresh for the Unicorn_Profits
# workbook: updates the cached market-price / profit columns (H-N) on
# specific leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR
# sheets to their newly-fetched values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1618.2
$ws.Range("I107").Value = 1742.4445
$ws.Range("K107").Value = 1742.4445
$ws.Range("M107").Value = 177.5554999999999

$ws.Range("H112").Value = 1707.45
$ws.Range("J112").Value = 2120
$ws.Range("L112").Value = 6360
$ws.Range("N112").Value = -8576

$ws.Range("H118").Value = 83867.914
$ws.Range("I118").Value = 83867.914
$ws.Range("K118").Value = 251603.742
$ws.Range("M118").Value = -249946.742

$ws.Range("H132").Value = 2320.0784
$ws.Range("I132").Value = 1345.8096
$ws.Range("J132").Value = 6866.6665
$ws.Range("K132").Value = 4037.4288
$ws.Range("L132").Value = 20599.9995
$ws.Range("M132").Value = -1507.4288
$ws.Range("N132").Value = -25659.9995

$ws.Range("H141").Value = 1634.1818
$ws.Range("I141").Value = 811.1786
$ws.Range("J141").Value = 6243
$ws.Range("K141").Value = 2433.5358
$ws.Range("L141").Value = 18729
$ws.Range("M141").Value = 2746.4642
$ws.Range("N141").Value = -29089

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1132702.9
$ws.Range("I32").Value = 1393534.6
$ws.Range("K32").Value = 1393534.6
$ws.Range("M32").Value = -1393247.6

$ws.Range("H61").Value = 422451.03
$ws.Range("I61").Value = 356564.84
$ws.Range("J61").Value = 529949.5600000001
$ws.Range("K61").Value = 356564.84
$ws.Range("L61").Value = 529949.5600000001
$ws.Range("M61").Value = -356352.84
$ws.Range("N61").Value = -530373.5600000001

$ws.Range("H74").Value = 170545.23
$ws.Range("I74").Value = 193092.81
$ws.Range("J74").Value = 80354.92
$ws.Range("K74").Value = 193092.81
$ws.Range("L74").Value = 80354.92
$ws.Range("M74").Value = -192218.81
$ws.Range("N74").Value = -82102.92

$ws.Range("H77").Value = 170545.23
$ws.Range("I77").Value = 193092.81
$ws.Range("J77").Value = 80354.92
$ws.Range("K77").Value = 965464.05
$ws.Range("L77").Value = 401774.6
$ws.Range("M77").Value = -961096.05
$ws.Range("N77").Value = -410510.6

$ws.Range("H88").Value = 4050
$ws.Range("I88").Value = 6250
$ws.Range("J88").Value = 2583.3333
$ws.Range("K88").Value = 6250
$ws.Range("L88").Value = 2583.3333
$ws.Range("M88").Value = -5844
$ws.Range("N88").Value = -3395.3333

$ws.Range("H91").Value = 4050
$ws.Range("I91").Value = 6250
$ws.Range("J91").Value = 2583.3333
$ws.Range("K91").Value = 6250
$ws.Range("L91").Value = 2583.3333
$ws.Range("M91").Value = -4846
$ws.Range("N91").Value = -5391.3333

$ws.Range("H132").Value = 19318.018
$ws.Range("I132").Value = 26021.072
$ws.Range("J132").Value = 3677.5557
$ws.Range("K132").Value = 78063.216
$ws.Range("L132").Value = 11032.6671
$ws.Range("M132").Value = -75533.216
$ws.Range("N132").Value = -16092.6671

$ws.Range("H136").Value = 422451.03
$ws.Range("I136").Value = 356564.84
$ws.Range("J136").Value = 529949.5600000001
$ws.Range("K136").Value = 1069694.52
$ws.Range("L136").Value = 1589848.68
$ws.Range("M136").Value = -1067144.52
$ws.Range("N136").Value = -1594948.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2602.7026
$ws.Range("I86").Value = 5709.091
$ws.Range("J86").Value = 1288.4615
$ws.Range("K86").Value = 5709.091
$ws.Range("L86").Value = 1288.4615
$ws.Range("M86").Value = -4586.091
$ws.Range("N86").Value = -3534.4615

$ws.Range("H89").Value = 2602.7026
$ws.Range("I89").Value = 5709.091
$ws.Range("J89").Value = 1288.4615
$ws.Range("K89").Value = 28545.455
$ws.Range("L89").Value = 6442.307499999999
$ws.Range("M89").Value = -22929.455
$ws.Range("N89").Value = -17674.3075

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2013.9048
$ws.Range("I132").Value = 966.9666999999999
$ws.Range("J132").Value = 4631.25
$ws.Range("K132").Value = 2900.9001
$ws.Range("L132").Value = 13893.75
$ws.Range("M132").Value = -370.9000999999998
$ws.Range("N132").Value = -18953.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 922.82355
$ws.Range("I103").Value = 269.5
$ws.Range("J103").Value = 1856.1428
$ws.Range("K103").Value = 808.5
$ws.Range("L103").Value = 5568.428400000001
$ws.Range("M103").Value = 70.5
$ws.Range("N103").Value = -7326.428400000001

$ws.Range("H131").Value = 1516.4423
$ws.Range("I131").Value = 1890
$ws.Range("J131").Value = 1458.3334
$ws.Range("K131").Value = 5670
$ws.Range("L131").Value = 4375.0002
$ws.Range("M131").Value = -630
$ws.Range("N131").Value = -14455.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4772.4893
$ws.Range("I70").Value = 3953.8076
$ws.Range("J70").Value = 5786.095
$ws.Range("K70").Value = 3953.8076
$ws.Range("L70").Value = 5786.095
$ws.Range("M70").Value = -3683.8076
$ws.Range("N70").Value = -6326.095

$ws.Range("H73").Value = 4772.4893
$ws.Range("I73").Value = 3953.8076
$ws.Range("J73").Value = 5786.095
$ws.Range("K73").Value = 3953.8076
$ws.Range("L73").Value = 5786.095
$ws.Range("M73").Value = -3017.8076
$ws.Range("N73").Value = -7658.095

$ws.Range("H97").Value = 875
$ws.Range("I97").Value = 800
$ws.Range("J97").Value = 950
$ws.Range("K97").Value = 800
$ws.Range("L97").Value = 950
$ws.Range("M97").Value = -304
$ws.Range("N97").Value = -1942

$ws.Range("H132").Value = 3704.075
$ws.Range("I132").Value = 3680.25
$ws.Range("J132").Value = 3759.6667
$ws.Range("K132").Value = 11040.75
$ws.Range("L132").Value = 11279.0001
$ws.Range("M132").Value = -8510.75
$ws.Range("N132").Value = -16339.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1891
$ws.Range("I40").Value = 1326.1666
$ws.Range("J40").Value = 5280
$ws.Range("K40").Value = 1326.1666
$ws.Range("L40").Value = 5280
$ws.Range("M40").Value = -1190.1666
$ws.Range("N40").Value = -5552

$ws.Range("H132").Value = 7967.1353
$ws.Range("I132").Value = 2376.6365
$ws.Range("J132").Value = 16166.533
$ws.Range("K132").Value = 7129.9095
$ws.Range("L132").Value = 48499.599
$ws.Range("M132").Value = -4599.9095
$ws.Range("N132").Value = -53559.599

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1362.0509
$ws.Range("I132").Value = 784.6667
$ws.Range("J132").Value = 3881.5454
$ws.Range("K132").Value = 2354.0001
$ws.Range("L132").Value = 11644.6362
$ws.Range("M132").Value = 175.9998999999998
$ws.Range("N132").Value = -16704.6362
